# Fruta / hortaliza, semanal
# Insert a new weekly record as row 58, pushing the existing rows 58-65
# down to 59-66 (dimension grows from A1:R65 to A1:R66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 58 - this shifts rows 58:65 down to 59:66.
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new weekly record.
$ws.Range("A58").Value = 2
$ws.Range("B58").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C58").Value = "Coquimbo"
$ws.Range("D58").Value = 44776
$ws.Range("E58").Value = 4
$ws.Range("F58").Value = 100112022
$ws.Range("G58").Value = "Arveja Verde"
$ws.Range("H58").Value = "Perfection"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 400
$ws.Range("K58").Value = 28000
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = 29000
$ws.Range("N58").Value = "`$/malla 25 kilos"
$ws.Range("O58").Value = "Provincia de Limarí"
$ws.Range("P58").Value = 1160
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"
